$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.379.87"
$ws.Range("E2").Value = "  -1.19%  "

$ws.Range("D3").Value = "2.047.85"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("D4").Value = "0.994"
$ws.Range("E4").Value = "  -0.66%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.00"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "56.84"
$ws.Range("E8").Value = "  -2.93%  "

$ws.Range("E9").Value = "  -2.24%  "

$ws.Range("D10").Value = "0.0811"
$ws.Range("E10").Value = "  +3.41%  "

$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.353.35"
$ws.Range("E12").Value = "  -1.66%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "14.59"
$ws.Range("E13").Value = "  -3.23%  "

$ws.Range("D14").Value = "20.69"
$ws.Range("E14").Value = "  -3.02%  "

$ws.Range("D15").Value = "0.756"
$ws.Range("E15").Value = "  -3.22%  "

$ws.Range("D16").Value = "5.29"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("D17").Value = "2.048.51"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "37.223.06"
$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("D19").Value = "6.02"
$ws.Range("E19").Value = "  -1.96%  "

$ws.Range("E20").Value = "  -1.78%  "

$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "226.66"
$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").Value = "2.37"
$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  -5.29%  "

$ws.Range("D26").Value = "9.51"
$ws.Range("E26").Value = "  -3.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.80"

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.129"
$ws.Range("E28").Value = "  -5.66%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "1.39"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").Value = "18.93"
$ws.Range("E30").Value = "  -3.19%  "

$ws.Range("E31").Value = "  -2.70%  "

$ws.Range("D32").Value = "4.53"
$ws.Range("E32").Value = "  -4.19%  "

$ws.Range("D33").Value = "4.59"
$ws.Range("E33").Value = "  -2.04%  "

$ws.Range("D34").Value = "0.0612"
$ws.Range("E34").Value = "  -3.48%  "

$ws.Range("E35").Value = "  -3.16%  "

$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("D37").Value = "0.996"
$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  -5.27%  "

$ws.Range("D39").Value = "5.38"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("D40").Value = "0.0221"
$ws.Range("E40").Value = "  -6.63%  "

$ws.Range("D41").Value = "17.02"
$ws.Range("E41").Value = "  +1.06%  "

$ws.Range("D42").Value = "1.482.37"
$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("D44").Value = "0.0942"
$ws.Range("E44").Value = "  -3.69%  "

$ws.Range("D45").Value = "96.26"
$ws.Range("E45").Value = "  -6.33%  "

$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("E47").Value = "  -4.50%  "

$ws.Range("D48").Value = "7.15"
$ws.Range("E48").Value = "  -2.46%  "

$ws.Range("D49").Value = "2.91"
$ws.Range("E49").Value = "  -2.98%  "

$ws.Range("D50").Value = "3.72"
$ws.Range("E50").Value = "  -9.81%  "

$ws.Range("D51").Value = "2.239.16"
$ws.Range("E51").Value = "  -1.62%  "
